$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 385
$ws.Range("G4").Value = 1497
$ws.Range("G5").Value = 2880
$ws.Range("G6").Value = 80006
$ws.Range("G7").Value = 60615
$ws.Range("G8").Value = 366126
$ws.Range("G9").Value = 579632
$ws.Range("G10").Value = 1749837

$ws.Range("G11").Select()
